$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text cells, plain text reassignment) ---
$ws.Range("A8").Value = "Volume 31   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/9/2024  Through  9/15/2024"

# --- Row 15 ---
$ws.Range("N15").Value = -44

# --- Row 16 ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -47.058823529411
$ws.Range("I16").Value = 82
$ws.Range("J16").Value = 104
$ws.Range("K16").Value = -21.153846153846
$ws.Range("L16").Value = -38.805970149253
$ws.Range("M16").Value = -56.14973262032
$ws.Range("N16").Value = -87.267080745341

# --- Row 17 ---
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 17
$ws.Range("H17").Value = -5.555555555555
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 158
$ws.Range("K17").Value = 26.582278481012
$ws.Range("L17").Value = 49.253731343283
$ws.Range("M17").Value = 203.030303030303
$ws.Range("N17").Value = 12.359550561797

# --- Row 18 ---
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -16.666666666666
$ws.Range("I18").Value = 154
$ws.Range("J18").Value = 196
$ws.Range("K18").Value = -21.428571428571
$ws.Range("L18").Value = -1.282051282051
$ws.Range("M18").Value = -18.518518518518
$ws.Range("N18").Value = -85.457979225684

# --- Row 19 ---
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -69.565217391304
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -34.426229508196
$ws.Range("I19").Value = 394
$ws.Range("J19").Value = 442
$ws.Range("K19").Value = -10.859728506787
$ws.Range("L19").Value = -16.525423728813
$ws.Range("M19").Value = 25.477707006369
$ws.Range("N19").Value = -7.294117647058

# --- Row 20 ---
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = 25.925925925925
$ws.Range("I20").Value = 254
$ws.Range("J20").Value = 260
$ws.Range("K20").Value = -2.307692307692
$ws.Range("L20").Value = 38.79781420765
$ws.Range("M20").Value = 67.105263157894
$ws.Range("N20").Value = -92.629135229251

# --- Row 21 ---
$ws.Range("D21").Value = 39
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 148
$ws.Range("H21").Value = -17.567567567567
$ws.Range("I21").Value = 1101
$ws.Range("J21").Value = 1172
$ws.Range("K21").Value = -6.058020477815
$ws.Range("L21").Value = 0.547945205479
$ws.Range("M21").Value = 19.673913043478
$ws.Range("N21").Value = -80.961438699636

# --- Row 22 ---
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 32
$ws.Range("K22").Value = -21.875
$ws.Range("L22").Value = 25

# --- Row 23 (C/D/E become blank-style text cells "0"/"0"/"***.*") ---
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E29").Copy($ws.Range("E23"))
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -40

# --- Row 24 ---
$ws.Range("C24").Value = 15
$ws.Range("E24").Value = -44.444444444444
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 126
$ws.Range("H24").Value = -32.539682539682
$ws.Range("I24").Value = 840
$ws.Range("J24").Value = 1021
$ws.Range("K24").Value = -17.727717923604
$ws.Range("L24").Value = -17.077986179664
$ws.Range("M24").Value = 19.148936170212

# --- Row 25 ---
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = -34.482758620689
$ws.Range("I25").Value = 359
$ws.Range("J25").Value = 378
$ws.Range("K25").Value = -5.026455026455
$ws.Range("L25").Value = 16.938110749185

# --- Row 26 ---
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 41
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 41.379310344827
$ws.Range("I26").Value = 362
$ws.Range("J26").Value = 330
$ws.Range("K26").Value = 9.696969696969
$ws.Range("L26").Value = 11.384615384615
$ws.Range("M26").Value = 25.694444444444

# --- Row 27 (D/E become blank-style text cells "0"/"***.*") ---
$ws.Range("C27").Copy($ws.Range("D27"))
$ws.Range("E29").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 0

# --- Row 28 (C/D/E convert from blank-style text back to numbers) ---
$ws.Range("F28").Copy($ws.Range("C28"))
$ws.Range("F28").Copy($ws.Range("D28"))
$ws.Range("H28").Copy($ws.Range("E28"))
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = -60
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = -12.121212121212
$ws.Range("L28").Value = -35.555555555555

# --- Row 29 ---
$ws.Range("L29").Value = -66.666666666666

# --- Row 30 ---
$ws.Range("L30").Value = -57.142857142857
